$d = $word.ActiveDocument

# Paragraph indices (1-based, Word COM convention) that currently read
# "Token nodig: Nee" and must become "Token nodig: Ja":
#   11 -> "GET rated movies" section
#   30 -> "GET users" section
#   36 -> "GET user" section
$targets = @(11, 30, 36)

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    # Trim the trailing paragraph mark from the range so Find only matches
    # inside the visible text.
    $r.End = $r.End - 1
    $r.Find.Execute("Token nodig: Nee", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "Token nodig: Ja", 2)
}

$d.Save()
